$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column: copy header formatting from the last existing header (G1)
# then overwrite the value so the style (bold, centered, bordered) matches.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data values for the Save column (plain numeric cells, no special style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
